# Apply "diagram labeling thru 2016" edit:
#  - append "(s)" to singular subcategory labels in column H
#    (photo -> photo(s), line graph -> line graph(s), bar chart -> bar chart(s), drawing -> drawing(s))
#  - remove the now-unused "is_viewed" column (I)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$pluralize = @{
    "photo"      = "photo(s)"
    "line graph" = "line graph(s)"
    "bar chart"  = "bar chart(s)"
    "drawing"    = "drawing(s)"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $val = $cell.Value2
    if ($pluralize.ContainsKey($val)) {
        $cell.Value = $pluralize[$val]
    }
}

# Remove column I (the "is_viewed" boolean flag column) entirely.
$ws.Columns.Item(9).Delete()
